# Auto-applies the "calapan and latest updates" edit to the Enduse List sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Enduse Name (column B) text for rows 3-255 to their final values ---
# (Several rows shift because new entries were inserted alphabetically into the
#  master Enduse list, and a couple of rows fix a duplicate-text bug.)
$ws.Cells.Item(3, 2).Value = "1.5 MVA Station Transformer DG4 & DG5 Generator Winding/VCB/MOCB/Cable Monitor"
$ws.Cells.Item(4, 2).Value = "125 Vdc Battery Charger"
$ws.Cells.Item(5, 2).Value = "24VDC Bank Batteries"
$ws.Cells.Item(6, 2).Value = "3 Units Exhaust Fan ( MCI )"
$ws.Cells.Item(7, 2).Value = "4 Units Pielstick Generator"
$ws.Cells.Item(8, 2).Value = "40 MVA Power Transformer"
$ws.Cells.Item(9, 2).Value = "40 MVA, 69kV/6.6kV Power Transformer"
$ws.Cells.Item(10, 2).Value = "69KV Pole No. 3 and 4"
$ws.Cells.Item(11, 2).Value = "750kVA Transformer"
$ws.Cells.Item(12, 2).Value = "750KVA Transformer Protection Relay"
$ws.Cells.Item(13, 2).Value = "Acetylene and Oxygen"
$ws.Cells.Item(14, 2).Value = "Adopt An Estero Clean-Up Activity on June 22, 2019"
$ws.Cells.Item(15, 2).Value = "Air Compressor"
$ws.Cells.Item(16, 2).Value = "Air Cooler"
$ws.Cells.Item(17, 2).Value = "Air Intake Room`n"
$ws.Cells.Item(18, 2).Value = "Allan Amoguis"
$ws.Cells.Item(19, 2).Value = "Asset Management"
$ws.Cells.Item(20, 2).Value = "Assorted Engine Parts"
$ws.Cells.Item(21, 2).Value = "Auxiliary Generator"
$ws.Cells.Item(22, 2).Value = "Auxiliary Lightings"
$ws.Cells.Item(23, 2).Value = "Bacolod Office"
$ws.Cells.Item(24, 2).Value = "Barring Gear Motor - Unit 2"
$ws.Cells.Item(25, 2).Value = "Barring Gear Motor - Unit 3"
$ws.Cells.Item(26, 2).Value = "Battery Charger"
$ws.Cells.Item(27, 2).Value = "Blood Sugar Monitoring and Clinic Supplies for Employee"
$ws.Cells.Item(28, 2).Value = "Boiler Circulating Pump and Motor Units 3 and 4"
$ws.Cells.Item(29, 2).Value = "Boiler Circulating Pump No.4"
$ws.Cells.Item(30, 2).Value = "Boiler Condensate Pump Motor"
$ws.Cells.Item(31, 2).Value = "Boiler Water Intake"
$ws.Cells.Item(32, 2).Value = "Borromeo's Lot"
$ws.Cells.Item(33, 2).Value = "Brigada Eskwela (Teodoro M. Morada Sr. Elem. School)"
$ws.Cells.Item(34, 2).Value = "Bus Differential Panel"
$ws.Cells.Item(35, 2).Value = "Bus Protection Relay"
$ws.Cells.Item(36, 2).Value = "Canteen`n"
$ws.Cells.Item(37, 2).Value = "Canteen and Ladies Dorm"
$ws.Cells.Item(38, 2).Value = "Cempco Office"
$ws.Cells.Item(39, 2).Value = "CENPRI Employees`n"
$ws.Cells.Item(40, 2).Value = "CENPRI Signage"
$ws.Cells.Item(41, 2).Value = "Cenpri Warehouse Building`n"
$ws.Cells.Item(42, 2).Value = "CENPRI Warehouse Department"
$ws.Cells.Item(43, 2).Value = "CENPRI Warehouse Extension"
$ws.Cells.Item(44, 2).Value = "CENPRI Warehouse Office"
$ws.Cells.Item(45, 2).Value = "Christmas Lantern"
$ws.Cells.Item(46, 2).Value = "Christmas Party Stage"
$ws.Cells.Item(47, 2).Value = "Circular Saw"
$ws.Cells.Item(48, 2).Value = "Clinic Use`n"
$ws.Cells.Item(49, 2).Value = "Common Bus Differential Fault & Breaker Failure"
$ws.Cells.Item(50, 2).Value = "Community"
$ws.Cells.Item(51, 2).Value = "Company Meeting (Visual Presentation)"
$ws.Cells.Item(52, 2).Value = "Computer/Electronic device power supply"
$ws.Cells.Item(53, 2).Value = "Conference Room"
$ws.Cells.Item(54, 2).Value = "Control Air Compressor - Common"
$ws.Cells.Item(55, 2).Value = "Control Air Compressor No. 1"
$ws.Cells.Item(56, 2).Value = "Control Air Compressor No. 2"
$ws.Cells.Item(57, 2).Value = "Control Air Compressor No. 3"
$ws.Cells.Item(58, 2).Value = "Control Panel"
$ws.Cells.Item(59, 2).Value = "Control Room"
$ws.Cells.Item(60, 2).Value = "Cooling Tower"
$ws.Cells.Item(61, 2).Value = "Cooling Tower Basin -  Common"
$ws.Cells.Item(62, 2).Value = "Cooling Tower Basin No. 1 "
$ws.Cells.Item(63, 2).Value = "Cooling Tower Basin No. 2"
$ws.Cells.Item(64, 2).Value = "Crane & Flat Bed Trailer"
$ws.Cells.Item(65, 2).Value = "Cummins Engine"
$ws.Cells.Item(66, 2).Value = "CV Access Bay Area"
$ws.Cells.Item(67, 2).Value = "Cylinder Head Assy & Cylinder At CV Access Bay Area"
$ws.Cells.Item(68, 2).Value = "Cylinder Head, Spare, Pielstick"
$ws.Cells.Item(69, 2).Value = "Deep Well Facility"
$ws.Cells.Item(70, 2).Value = "Deep Well Pump"
$ws.Cells.Item(71, 2).Value = "Deep Well Riser Pipes Pull-out"
$ws.Cells.Item(72, 2).Value = "DG 1 and 2"
$ws.Cells.Item(73, 2).Value = "DG Unit 4 Linkage Clamp"
$ws.Cells.Item(74, 2).Value = "DG1 `n"
$ws.Cells.Item(75, 2).Value = "DG1 (CV Area) Main Engine Parts & Components"
$ws.Cells.Item(76, 2).Value = "DG2`n"
$ws.Cells.Item(77, 2).Value = "DG3"
$ws.Cells.Item(78, 2).Value = "DG4"
$ws.Cells.Item(79, 2).Value = "DG4(CV Area) Main Engine Parts & Components"
$ws.Cells.Item(80, 2).Value = "DG5"
$ws.Cells.Item(81, 2).Value = "Diesel Storage Tank"
$ws.Cells.Item(82, 2).Value = "Distillation Equipment"
$ws.Cells.Item(83, 2).Value = "Drain Pipe in Running Units Sulzer and Smoke Stack"
$ws.Cells.Item(84, 2).Value = "Drum Table & Chairs"
$ws.Cells.Item(85, 2).Value = "Electrical Consumables"
$ws.Cells.Item(86, 2).Value = "Electrical Equipment"
$ws.Cells.Item(87, 2).Value = "Electrical Handtools"
$ws.Cells.Item(88, 2).Value = "Electrical Instruments and Plant Equipment Protection"
$ws.Cells.Item(89, 2).Value = "Engine Auxiliary Area Lighting"
$ws.Cells.Item(90, 2).Value = "Engine Auxiliary Lgihting & Power Panel"
$ws.Cells.Item(91, 2).Value = "Engine Auxiliary Lighting & Power Panel"
$ws.Cells.Item(92, 2).Value = "Engine Driven"
$ws.Cells.Item(93, 2).Value = "Environment and Pollution Control"
$ws.Cells.Item(94, 2).Value = "Fabrication of Powerhouse Ventilation Louvers"
$ws.Cells.Item(95, 2).Value = "Facilities Improvement"
$ws.Cells.Item(96, 2).Value = "Fire Brigade Training"
$ws.Cells.Item(97, 2).Value = "Fire Fighting System"
$ws.Cells.Item(98, 2).Value = "Fire Hose Cabinet/Dead End of Fire Sprinkler of Piping"
$ws.Cells.Item(99, 2).Value = "Fire Hydrant"
$ws.Cells.Item(100, 2).Value = "Fire Safety"
$ws.Cells.Item(101, 2).Value = "Flag Pole"
$ws.Cells.Item(102, 2).Value = "Fluke Clamp Meter"
$ws.Cells.Item(103, 2).Value = "Fuel and Lube Oil Management"
$ws.Cells.Item(104, 2).Value = "Fuel and Lube Oil Recovery System"
$ws.Cells.Item(105, 2).Value = "Fuel Farm"
$ws.Cells.Item(106, 2).Value = "Fuel Module 1 and 3, Starting Air Compressor Unit 3 and 5"
$ws.Cells.Item(107, 2).Value = "Fuel Module Booster Pump Unit 2"
$ws.Cells.Item(108, 2).Value = "Fuel Recovery & Sludge Tanks"
$ws.Cells.Item(109, 2).Value = "Fuel Tank"
$ws.Cells.Item(110, 2).Value = "Garbage Bin"
$ws.Cells.Item(111, 2).Value = "Generating Unit Protection & EIC Equipment"
$ws.Cells.Item(112, 2).Value = "Generator Master Panel"
$ws.Cells.Item(113, 2).Value = "Generator Sliding Tools"
$ws.Cells.Item(114, 2).Value = "Generator Unit 1"
$ws.Cells.Item(115, 2).Value = "Generator Unit 2"
$ws.Cells.Item(116, 2).Value = "Generator Unit 4"
$ws.Cells.Item(117, 2).Value = "Grounding System Lay-out / Installation"
$ws.Cells.Item(118, 2).Value = "Guardhouse`n"
$ws.Cells.Item(119, 2).Value = "Guests/VIP"
$ws.Cells.Item(120, 2).Value = "Heavy Equipment - Boomtruck"
$ws.Cells.Item(121, 2).Value = "Heavy Equipment - Boomtruck and Forklift"
$ws.Cells.Item(122, 2).Value = "Heavy Equipment - Forklift"
$ws.Cells.Item(123, 2).Value = "Heavy Fuel Oil Purifier - Common"
$ws.Cells.Item(124, 2).Value = "Heavy Fuel Oil Purifier No. 1"
$ws.Cells.Item(125, 2).Value = "Heavy Fuel Oil Purifier No. 2"
$ws.Cells.Item(126, 2).Value = "Heavy Fuel Oil Sludge Basin"
$ws.Cells.Item(127, 2).Value = "HFO Circulating Pump and Motor Coupling"
$ws.Cells.Item(128, 2).Value = "HFO Fuel Piping Insulation and Cladding"
$ws.Cells.Item(129, 2).Value = "HFO Recovery Tank"
$ws.Cells.Item(130, 2).Value = "HFO Settling & Service Tanks"
$ws.Cells.Item(131, 2).Value = "Honing Machine"
$ws.Cells.Item(132, 2).Value = "Honing Machine Compressor Motor"
$ws.Cells.Item(133, 2).Value = "Isuzu, Pick-up, Fuego, Diesel, 1999"
$ws.Cells.Item(134, 2).Value = "Jacket Water Cooler - Common"
$ws.Cells.Item(135, 2).Value = "Jacket Water Cooler No. 1"
$ws.Cells.Item(136, 2).Value = "Jacket Water Cooler No. 2"
$ws.Cells.Item(137, 2).Value = "Jacket Water Cooler No. 3"
$ws.Cells.Item(138, 2).Value = "Jacket Water Cooler No. 4"
$ws.Cells.Item(139, 2).Value = "Jacket Water Cooler No. 5"
$ws.Cells.Item(140, 2).Value = "Jacket Water Motor #4, MCCB"
$ws.Cells.Item(141, 2).Value = "Jacket Water Pump Motor No.4"
$ws.Cells.Item(142, 2).Value = "Laboratory Use"
$ws.Cells.Item(143, 2).Value = "Ladies' Dormitory`n"
$ws.Cells.Item(144, 2).Value = "Lifting Equipment for Transferring of Heavy Parts/Boxes with Parts"
$ws.Cells.Item(145, 2).Value = "Lube Oil Cooler - Common"
$ws.Cells.Item(146, 2).Value = "Lube Oil Cooler No. 1"
$ws.Cells.Item(147, 2).Value = "Lube Oil Cooler No. 2"
$ws.Cells.Item(148, 2).Value = "Lube Oil Cooler No. 3"
$ws.Cells.Item(149, 2).Value = "Lube Oil Cooler No. 4"
$ws.Cells.Item(150, 2).Value = "Lube Oil Cooler No. 5"
$ws.Cells.Item(151, 2).Value = "Lube Oil Priming Pump"
$ws.Cells.Item(152, 2).Value = "Lube Oil Priming Pump"
$ws.Cells.Item(153, 2).Value = "Lube Oil Purifier - Common"
$ws.Cells.Item(154, 2).Value = "Lube Oil Purifier No. 1"
$ws.Cells.Item(155, 2).Value = "Lube Oil Purifier No. 2"
$ws.Cells.Item(156, 2).Value = "Lube Oil Purifier No. 3"
$ws.Cells.Item(157, 2).Value = "Lube Oil Purifier No. 4"
$ws.Cells.Item(158, 2).Value = "Lube Oil Purifier No. 5"
$ws.Cells.Item(159, 2).Value = "Maintenance and Operation"
$ws.Cells.Item(160, 2).Value = "Maintenance Reconditioning Area and Fuel Farm Area"
$ws.Cells.Item(161, 2).Value = "Maintenance Tools"
$ws.Cells.Item(162, 2).Value = "Male Common CR"
$ws.Cells.Item(163, 2).Value = "Master Control Panel"
$ws.Cells.Item(164, 2).Value = "Mechanical Barracks"
$ws.Cells.Item(165, 2).Value = "Microwave Antenna"
$ws.Cells.Item(166, 2).Value = "Microwave Control Panel"
$ws.Cells.Item(167, 2).Value = "Microwave Panel"
$ws.Cells.Item(168, 2).Value = "Microwave Radio Equipment"
$ws.Cells.Item(169, 2).Value = "MOCB Units 4 & 5"
$ws.Cells.Item(170, 2).Value = "Mono Pump Chamber No.1"
$ws.Cells.Item(171, 2).Value = "Motor Control Center 3"
$ws.Cells.Item(172, 2).Value = "NALCO and SEM Water Softener"
$ws.Cells.Item(173, 2).Value = "NALCO Water Softener Unit"
$ws.Cells.Item(174, 2).Value = "New 750kVA Station Service Transformer"
$ws.Cells.Item(175, 2).Value = "Non Disturbance Monitoring Equipment"
$ws.Cells.Item(176, 2).Value = "NVR CCTV Cameras"
$ws.Cells.Item(177, 2).Value = "Office Use`n"
$ws.Cells.Item(178, 2).Value = "Oil Analysis Equipment"
$ws.Cells.Item(179, 2).Value = "Oil and Water Mechanical Separator Basin"
$ws.Cells.Item(180, 2).Value = "Operations & Maintenance Consumables"
$ws.Cells.Item(181, 2).Value = "Operations Communications"
$ws.Cells.Item(182, 2).Value = "Out-of-House Oil Analysis"
$ws.Cells.Item(183, 2).Value = "Panel Main Breaker"
$ws.Cells.Item(184, 2).Value = "Pedestrian Crossing In-front of Main Gate"
$ws.Cells.Item(185, 2).Value = "Pielstick Jacket Liner"
$ws.Cells.Item(186, 2).Value = "Plant Common Tools, Special Tools and Equipment"
$ws.Cells.Item(187, 2).Value = "Plant Common Tools, Special Tools, and Equipment"
$ws.Cells.Item(188, 2).Value = "Plant Decoration"
$ws.Cells.Item(189, 2).Value = "Plant Equipment"
$ws.Cells.Item(190, 2).Value = "Plant Site Security & Monitoring Equipment"
$ws.Cells.Item(191, 2).Value = "Plant Testing Instrument"
$ws.Cells.Item(192, 2).Value = "Plate Compactor"
$ws.Cells.Item(193, 2).Value = "Power Plant Premises"
$ws.Cells.Item(194, 2).Value = "Powerhouse -  Auxiliary Side"
$ws.Cells.Item(195, 2).Value = "Powerhouse Building`n"
$ws.Cells.Item(196, 2).Value = "Powerhouse Toolbox"
$ws.Cells.Item(197, 2).Value = "Progen Office Use"
$ws.Cells.Item(198, 2).Value = "Progen Warehouse`n"
$ws.Cells.Item(199, 2).Value = "Raw Water Unit 4 & 5"
$ws.Cells.Item(200, 2).Value = "Reconditioning Equipment"
$ws.Cells.Item(201, 2).Value = "Reconditioning Precision Tools"
$ws.Cells.Item(202, 2).Value = "Recovered Materials"
$ws.Cells.Item(203, 2).Value = "Relief Valve Cap for Cylinder Head"
$ws.Cells.Item(204, 2).Value = "Repainting of Color Coded Waste Bin"
$ws.Cells.Item(205, 2).Value = "Restrooms`n"
$ws.Cells.Item(206, 2).Value = "Running Units`n"
$ws.Cells.Item(207, 2).Value = "Running Units - Pielstick`n"
$ws.Cells.Item(208, 2).Value = "Running Units - Sulzer`n"
$ws.Cells.Item(209, 2).Value = "Running Units / Unit 1 Servicing Works"
$ws.Cells.Item(210, 2).Value = "Running Units Maintenance Tools"
$ws.Cells.Item(211, 2).Value = "Safety"
$ws.Cells.Item(212, 2).Value = "SCADA Room"
$ws.Cells.Item(213, 2).Value = "SEM Water Softener Unit"
$ws.Cells.Item(214, 2).Value = "Service Vehicle Isuzu Crosswind Plate No. FFN706"
$ws.Cells.Item(215, 2).Value = "Settling and Service Tank"
$ws.Cells.Item(216, 2).Value = "Sludge Tank"
$ws.Cells.Item(217, 2).Value = "Smoke Stack"
$ws.Cells.Item(218, 2).Value = "Soft Water Supply Pump"
$ws.Cells.Item(219, 2).Value = "Spare Ideal Generator"
$ws.Cells.Item(220, 2).Value = "Spare Stator"
$ws.Cells.Item(221, 2).Value = "Spare Stator Rewinding Enclosure"
$ws.Cells.Item(222, 2).Value = "Staffhouse 1 - Site"
$ws.Cells.Item(223, 2).Value = "Staffhouse 2 - Jara"
$ws.Cells.Item(224, 2).Value = "Staffhouse 3 - EDJ"
$ws.Cells.Item(225, 2).Value = "Starting Air Compressor - Common"
$ws.Cells.Item(226, 2).Value = "Starting Air Compressor - Pielstick"
$ws.Cells.Item(227, 2).Value = "Starting Air Compressor No. 1"
$ws.Cells.Item(228, 2).Value = "Starting Air Compressor No. 2"
$ws.Cells.Item(229, 2).Value = "Starting Air Compressor No. 3"
$ws.Cells.Item(230, 2).Value = "Starting Air Compressor No. 4"
$ws.Cells.Item(231, 2).Value = "Starting Air Compressor No. 5"
$ws.Cells.Item(232, 2).Value = "Station Load Metering"
$ws.Cells.Item(233, 2).Value = "Station Transformer Protection Relay"
$ws.Cells.Item(234, 2).Value = "Steam Equipment"
$ws.Cells.Item(235, 2).Value = "Substation`n"
$ws.Cells.Item(236, 2).Value = "Switch Gear Room"
$ws.Cells.Item(237, 2).Value = "Switch Yard"
$ws.Cells.Item(238, 2).Value = "Synchronizing Panel Units"
$ws.Cells.Item(239, 2).Value = "Tagging of Tools and Equipments"
$ws.Cells.Item(240, 2).Value = "Tank Farm`n"
$ws.Cells.Item(241, 2).Value = "Tank Farm Earth Grounding"
$ws.Cells.Item(242, 2).Value = "Testing Equipment Storage"
$ws.Cells.Item(243, 2).Value = "Testing of 750KVA Transformer"
$ws.Cells.Item(244, 2).Value = "Tools, Inventory-Mary Grace Bugna"
$ws.Cells.Item(245, 2).Value = "Toyota Inova"
$ws.Cells.Item(246, 2).Value = "Trading"
$ws.Cells.Item(247, 2).Value = "Trainees on Oil Spill"
$ws.Cells.Item(248, 2).Value = "Turbo Charger Air Intake"
$ws.Cells.Item(249, 2).Value = "Turning Gear Motor DG 1-3"
$ws.Cells.Item(250, 2).Value = "Unit 1 - 5 Synchronizing Panels"
$ws.Cells.Item(251, 2).Value = "Unit 5 Control Air System"
$ws.Cells.Item(252, 2).Value = "Unit Panel"
$ws.Cells.Item(253, 2).Value = "Warehouse Beginning Balance"
$ws.Cells.Item(254, 2).Value = "Warehouse Building"
$ws.Cells.Item(255, 2).Value = "Warehouse Extension Gate"

# --- Append 6 new rows (256-261) for entries that no longer fit before row 255 ---
$ws.Range("A255:B255").Copy($ws.Range("A256:B261"))
$ws.Cells.Item(256, 1).Value = 254
$ws.Cells.Item(256, 2).Value = "Waste Heat Recovery Boiler - Common"
$ws.Cells.Item(257, 1).Value = 255
$ws.Cells.Item(257, 2).Value = "Waste Heat Recovery Boiler No.1"
$ws.Cells.Item(258, 1).Value = 256
$ws.Cells.Item(258, 2).Value = "Waste Heat Recovery Boiler No.2"
$ws.Cells.Item(259, 1).Value = 257
$ws.Cells.Item(259, 2).Value = "Waste Heat Recovery Boiler No.3"
$ws.Cells.Item(260, 1).Value = 258
$ws.Cells.Item(260, 2).Value = "Westfalia Separator AG-Lube Oil"
$ws.Cells.Item(261, 1).Value = 259
$ws.Cells.Item(261, 2).Value = "Wire Marker Device"

# --- Match the saved selection state from the source workbook ---
$ws.Range("B261").Select()
